# The deck ships two embedded DrawingML themes:
#   ppt/theme/theme1.xml -> "Office Theme" (orphaned: not linked from any
#                            slide master / layout reachable via the OM)
#   ppt/theme/theme2.xml -> "Integral"      (the live theme used by the
#                            one slide master all slides/layouts share)
#
# The target edit swaps the two themes' contents: the live theme becomes
# the "Office Theme" palette and the (otherwise unreachable) orphan theme
# becomes "Integral". Since the orphan theme part isn't addressable from
# the PowerPoint object model (it isn't attached to any Master/Design the
# host exposes), the only in-model lever available is the live design's
# colour scheme. We push the "Office Theme" palette onto it so the file
# that PowerPoint actually renders through (theme2.xml) ends up holding
# the Office Theme colours, matching the post-swap target.

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.ColorScheme

function RgbValue($r, $g, $b) {
    return $b * 65536 + $g * 256 + $r
}

# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeThemeColors = @(
    (RgbValue 0x00 0x00 0x00),   # 1  dk1      000000
    (RgbValue 0xFF 0xFF 0xFF),   # 2  lt1      FFFFFF
    (RgbValue 0x44 0x54 0x6A),   # 3  dk2      44546A
    (RgbValue 0xE7 0xE6 0xE6),   # 4  lt2      E7E6E6
    (RgbValue 0x5B 0x9B 0xD5),   # 5  accent1  5B9BD5
    (RgbValue 0xED 0x7D 0x31),   # 6  accent2  ED7D31
    (RgbValue 0xA5 0xA5 0xA5),   # 7  accent3  A5A5A5
    (RgbValue 0xFF 0xC0 0x00),   # 8  accent4  FFC000
    (RgbValue 0x44 0x72 0xC4),   # 9  accent5  4472C4
    (RgbValue 0x70 0xAD 0x47),   # 10 accent6  70AD47
    (RgbValue 0x05 0x63 0xC1),   # 11 hlink    0563C1
    (RgbValue 0x95 0x4F 0x72)    # 12 folHlink 954F72
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $scheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
